# Rename "Sheet4" to "final" (the only substantive content change in the
# target diff; the remaining differences -- fileVersion/rupBuild bump,
# AlternateContent/x15ac absPath block, bookViews window geometry, calcPr
# concurrentCalc, extLst workbookPr/ArchID, fills-list pruning, theme font
# substitutions, and x14ac:dyDescent/sheetFormatPr tweaks -- are artifacts
# of the workbook being re-saved by a different Excel build and are not
# data/content edits.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Name = "final"
